$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Professional summary paragraph: neutralize "all Black and Asian-American"
#    voters language to "50M voters".
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML algorithms",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML algorithms", 2)

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: same language change, but "50M" must land in its
#    own bold / colored run (matching the other stat call-outs in the bullet).
#    First do a plain text swap, then re-split & format the "50M" run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial machine learning",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial machine learning", 2)

$numRange = $d.Content
$null = $numRange.Find.Execute("50M voters, developed geospatial machine")
$numRange.End = $numRange.Start + 3
$numRange.Font.Bold = 1
$numRange.Font.Color = 5258796

# ---------------------------------------------------------------------------
# 3) Move the "Field Director - The Feldman Group" experience block (heading +
#    4 paragraphs) from just before "KEY PROJECTS" to just before
#    "Software Engineer - Salsa Labs".
# ---------------------------------------------------------------------------
$srcStart = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "^Field Director - The Feldman Group") {
        $srcStart = $i
    }
}
$srcEnd = $srcStart + 4

$blockStart = $d.Paragraphs($srcStart).Range.Start
$blockEnd = $d.Paragraphs($srcEnd).Range.End
$blockRange = $d.Range($blockStart, $blockEnd)
$blockRange.Copy()

$destIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "Software Engineer - Salsa Labs") {
        $destIdx = $i
    }
}

$insertAt = $d.Paragraphs($destIdx).Range.Start
$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.Paste()

# The pasted heading paragraph loses its Heading 3 style on insertion; restore it.
$d.Paragraphs($destIdx).Style = "Heading 3"

# Remove the original copy of the block, which has shifted down by 5 paragraphs.
$oldStart = $srcStart + 5
$oldEnd = $srcEnd + 5
$oldBlockStart = $d.Paragraphs($oldStart).Range.Start
$oldBlockEnd = $d.Paragraphs($oldEnd).Range.End
$d.Range($oldBlockStart, $oldBlockEnd).Delete()

# ---------------------------------------------------------------------------
# 4) Key Projects "Geospatial Demographic Classification System" impact line.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved", 2)

Write-Output "done"
